$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header label cell A4 changes from "No." to "Num"
$ws.Range("A4").Value = "Num"

# 2. The "Matrial Description" header cell (C4) gets a Text number format
#    and left horizontal alignment (in addition to its existing
#    vertical-center + wrap formatting).
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").HorizontalAlignment = -4131

# 3. Build new info row 19 as a copy of row 1 (cell by cell, so existing
#    style indexes are reused instead of duplicated), then overwrite the
#    two label cells with their new text.
for ($col = 1; $col -le 9; $col++) {
  $src = $ws.Cells.Item(1, $col)
  $dst = $ws.Cells.Item(19, $col)
  $src.Copy($dst)
}
$ws.Range("A19").Value = "Company Code"
$ws.Range("E19").Value = "Test"

# 4. Build new info row 20 the same way, based on row 1 as well (row 1 and
#    row 2 share identical formatting).
for ($col = 1; $col -le 9; $col++) {
  $src = $ws.Cells.Item(1, $col)
  $dst = $ws.Cells.Item(20, $col)
  $src.Copy($dst)
}
$ws.Range("A20").Value = "Company"
$ws.Range("E20").Value = "Test1"

# 5. Re-create the merged cells for the two new info rows.
$ws.Range("A19:B19").MergeCells = $true
$ws.Range("C19:D19").MergeCells = $true
$ws.Range("E19:F19").MergeCells = $true
$ws.Range("G19:I19").MergeCells = $true

$ws.Range("A20:B20").MergeCells = $true
$ws.Range("C20:D20").MergeCells = $true
$ws.Range("E20:F20").MergeCells = $true
$ws.Range("G20:I20").MergeCells = $true

# 6. Row 21 is a duplicate of the (already updated) column-header row 4.
for ($col = 1; $col -le 9; $col++) {
  $src = $ws.Cells.Item(4, $col)
  $dst = $ws.Cells.Item(21, $col)
  $src.Copy($dst)
}

# 7. Rows 22 and 23 duplicate the data-row style/content used by rows 5-18.
for ($col = 1; $col -le 9; $col++) {
  $src = $ws.Cells.Item(5, $col)
  $dst = $ws.Cells.Item(22, $col)
  $src.Copy($dst)
}
for ($col = 1; $col -le 9; $col++) {
  $src = $ws.Cells.Item(5, $col)
  $dst = $ws.Cells.Item(23, $col)
  $src.Copy($dst)
}

# 8. Update the selection to match the new active cell.
[void]$ws.Range("A21").Select()
